$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.983.06"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "3.410.50"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'404.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "'131.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "'0.671"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").Value = "'0.121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("D11").Value = "'42.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "3.922.38"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "'8.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "'19.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "3.398.47"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "61.875.81"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "'10.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  -6.86%  "
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("D22").Value = "'84.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").Value = "'316.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'12.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("D25").Value = "'3.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "'4.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.02%  "
$ws.Range("D27").Value = "'29.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'7.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'2.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("D32").Value = "'0.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'41.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("D36").Value = "'0.0480"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").Value = "'51.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'3.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").Value = "'2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").Value = "'138.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").Value = "'1.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "'0.292"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "'3.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "'16.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "'21.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").Value = "2.125.72"
$ws.Range("E49").Value = "  -4.56%  "
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("D51").Value = "'1.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "

Write-Host "Applied crypto price/volume updates"
